{"js": "// Split the single \"A bibliografia...\" run into the intro sentence followed\n// by a blank line and five numbered references, each separated by a line\n// break (<w:br/>) within the SAME run - matching the target OOXML which\n// keeps everything inside one <w:r> but alternates <w:t>/<w:br/> children.\n//\n// Word's Office.js text model represents a line break (<w:br/>) as the\n// vertical-tab character U+000B inside a Range/Paragraph \".text\" string, and\n// Range.insertText(...) round-trips that back into <w:t>/<w:br/> runs, so we\n// build the replacement string with \"\\u000b\" at every point the diff adds a\n// <w:br/>.\n\nconst LBREAK = \"\\u000b\";\n\nconst introText =\n  \"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis e obtida na busca realizada pelos pr\u00f3prios alunos no in\u00edcio dos projetos. Seguem refer\u00eancias no t\u00f3pico de mentoria.\";\n\nconst refs = [\n  \"[1] Peddy, S. The art of mentoring \u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\",\n  \"[2] Zachary, L. J. The Mentor\\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\",\n  \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\",\n  \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\",\n  \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\",\n];\n\nconst newText =\n  introText + LBREAK + LBREAK + refs.join(LBREAK);\n\n// Locate the target paragraph by searching for the (unique) opening of its\n// text, then replace its whole range content in one shot so the run keeps\n// a single <w:r> wrapping every <w:t>/<w:br/> pair, exactly like the diff.\nconst results = context.document.body.search(\"A bibliografia ser\u00e1 recomendada\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst paragraph = results.items[0].paragraphs.getFirst();\nconst range = paragraph.getRange();\nrange.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split the single \"A bibliografia...\" run into the intro sentence followed\n# by a blank line and five numbered references, each separated by a manual\n# line break (<w:br/>) - all still inside the SAME run, matching the target\n# OOXML which keeps everything inside one <w:r> but alternates <w:t>/<w:br/>\n# children.\n#\n# In the Word object model a manual line break is Chr(11) (vertical tab /\n# vbVerticalTab); assigning a Range's .Text with that character embedded\n# serialises back to <w:t>.../ w:t><w:br/><w:t>...</w:t> runs, so we build the\n# replacement string with that character at every point the diff adds a\n# <w:br/>.\n\n$d = $word.ActiveDocument\n\n$lb = [char]11\n\n$introText = \"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis e obtida na busca realizada pelos pr\u00f3prios alunos no in\u00edcio dos projetos. Seguem refer\u00eancias no t\u00f3pico de mentoria.\"\n\n$ref1 = \"[1] Peddy, S. The art of mentoring \u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\"\n$ref2 = \"[2] Zachary, L. J. The Mentor\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\"\n$ref3 = \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\"\n$ref4 = \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\"\n$ref5 = \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\"\n\n$newText = $introText + $lb + $lb + $ref1 + $lb + $ref2 + $lb + $ref3 + $lb + $ref4 + $lb + $ref5\n\n# Locate the target paragraph by its (unique) opening text rather than a\n# hard-coded index, then overwrite its whole range in one shot so the\n# run keeps a single <w:r> wrapping every <w:t>/<w:br/> pair.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"A bibliografia ser\u00e1 recomendada*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Target paragraph not found\"\n}\n\n$target.Range.Text = $newText\n"}
